$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - add CNPJ/CPF Empresa (col K) and fix Descricao Portaria Saida (col M)
$ws.Range("K2").Value = "02039741070"
$ws.Range("M2").Value = "Portaria 2"

# Row 3 - add CNPJ/CPF Empresa (col K) and fix Descricao Portaria (cols M,N)
$ws.Range("K3").Value = "44456465654465"
$ws.Range("M3").Value = "Portaria 2"
$ws.Range("N3").Value = "Portaria 2"

# Row 4 - this was incomplete test data; clear the extra columns
$ws.Range("C4").Clear()
$ws.Range("D4").Clear()
$ws.Range("F4").Clear()
$ws.Range("G4").Clear()
$ws.Range("H4").Clear()
$ws.Range("I4").Clear()
$ws.Range("J4").Clear()
$ws.Range("N4").Clear()

# Row 5 - fix Data Saida, Hora Saida, add CNPJ/CPF + Nome Empresa, fix Descricao Portaria Saida
$ws.Range("C5").Value = "29/10/2023"
$ws.Range("D5").Value = "09:04"
$ws.Range("K5").Value = "17834987361926"
$ws.Range("L5").Value = "empresa teste"
$ws.Range("N5").Value = "Portaria 3"

# Row 6 - clear the extra columns
$ws.Range("C6").Clear()
$ws.Range("D6").Clear()
$ws.Range("F6").Clear()
$ws.Range("G6").Clear()
$ws.Range("H6").Clear()
$ws.Range("I6").Clear()
$ws.Range("N6").Clear()

# Row 7 - add CNPJ/CPF Empresa + Nome Empresa
$ws.Range("K7").Value = "28937465611"
$ws.Range("L7").Value = "Empresa testando"

# Row 8 - add CNPJ/CPF Empresa + Nome Empresa, add Descricao Portaria Saida
$ws.Range("K8").Value = "19283764532761"
$ws.Range("L8").Value = "Import test"
$ws.Range("N8").Value = "Portaria 1"

# Row 9 - clear the extra columns
$ws.Range("C9").Clear()
$ws.Range("D9").Clear()
$ws.Range("F9").Clear()
$ws.Range("G9").Clear()
$ws.Range("H9").Clear()
$ws.Range("I9").Clear()
$ws.Range("J9").Clear()
$ws.Range("N9").Clear()

# Update selection to match the author's final cursor position
$ws.Range("L9").Select()
